$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text in cell E8 and select it (mirrors user editing it in Excel)
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
